# Add files via upload
# Updates CO2Cap (sheet2) and CO2Price (sheet3) with an extra data row (period 8),
# and revises the AvailableBioEnergy (sheet4) figures plus appends its own period 8 row.

$wb = $excel.ActiveWorkbook

# --- CO2Cap sheet: append period 8 row ---
$wsCO2Cap = $wb.Worksheets.Item("CO2Cap")
$wsCO2Cap.Range("A11").Value = 8
$wsCO2Cap.Range("B11").Value = 213.273381294964

# --- CO2Price sheet: append period 8 row ---
$wsCO2Price = $wb.Worksheets.Item("CO2Price")
$wsCO2Price.Range("A11").Value = 8
$wsCO2Price.Range("B11").Value = 304.7619047619048

# --- AvailableBioEnergy sheet: revise existing values and append period 8 row ---
$wsBioEnergy = $wb.Worksheets.Item("AvailableBioEnergy")
$wsBioEnergy.Range("B4").Value = 2280000000
$wsBioEnergy.Range("B5").Value = 2368000000
$wsBioEnergy.Range("B6").Value = 2500000000
$wsBioEnergy.Range("B7").Value = 2500000000
$wsBioEnergy.Range("B8").Value = 2528000000
$wsBioEnergy.Range("B9").Value = 2612000000
$wsBioEnergy.Range("B10").Value = 2640000000
$wsBioEnergy.Range("A11").Value = 8
$wsBioEnergy.Range("B11").Value = 2640000000
